$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01384939437003609
$ws.Range("C2").Value = 2.218786637292586
$ws.Range("D2").Value = 21.84152128370215
$ws.Range("E2").Value = 4.673491337715536
$ws.Range("F2").Value = 4.783449883342865
$ws.Range("G2").Value = 22
$ws.Range("B3").Value = -0.1108952796442341
$ws.Range("C3").Value = 2.528008180603066
$ws.Range("D3").Value = 20.81929152853822
$ws.Range("E3").Value = 4.562816183952431
$ws.Range("F3").Value = 4.674114189209051
$ws.Range("G3").Value = 21
$ws.Range("B4").Value = -0.5909554399684317
$ws.Range("C4").Value = 1.803239635867731
$ws.Range("D4").Value = 9.880027863097968
$ws.Range("E4").Value = 3.143251161313389
$ws.Range("F4").Value = 3.167399652557398
$ws.Range("G4").Value = 20
$ws.Range("B5").Value = -0.1812640771482746
$ws.Range("C5").Value = 1.430345698091154
$ws.Range("D5").Value = 9.69518460969652
$ws.Range("E5").Value = 3.113709140189
$ws.Range("F5").Value = 3.19360672922055
$ws.Range("G5").Value = 19
$ws.Range("B6").Value = -0.1687012678555518
$ws.Range("C6").Value = 1.931467348941691
$ws.Range("D6").Value = 13.19207313971304
$ws.Range("E6").Value = 3.632089362847924
$ws.Range("F6").Value = 3.733355487989051
$ws.Range("G6").Value = 18
$ws.Range("B7").Value = -0.2974797958422067
$ws.Range("C7").Value = 1.914100007500023
$ws.Range("D7").Value = 11.03882254836823
$ws.Range("E7").Value = 3.322472354793675
$ws.Range("F7").Value = 3.410971099173742
$ws.Range("G7").Value = 17
$ws.Range("B8").Value = -0.1960625496366059
$ws.Range("C8").Value = 1.836194132777806
$ws.Range("D8").Value = 11.49911931585798
$ws.Range("E8").Value = 3.391035139283871
$ws.Range("F8").Value = 3.496387284801343
$ws.Range("G8").Value = 16